# "fixed happy planet index"
#
# 1. The "SP.POP.TOTL  Population, total," bullet (under POPULATION) was
#    marked up in red, like the "SM.POP.TOTL.ZS" bullet right below it.
# 2. The stray "_GoBack" bookmark (previously sitting on its own empty
#    paragraph after the SENSE OF SECURITY bullets) was moved into the
#    middle of the "Population density (people per sq. km of land area)"
#    bullet, right after the words "Population density".

$d = $word.ActiveDocument

# --- 1. Colour the "SP.POP.TOTL" bullet red -------------------------------
$rng = $d.Range(0, $d.Content.End)
$found = $rng.Find.Execute("SP.POP.TOTL")
if ($found) {
    $para = $rng.Paragraphs(1)
    $para.Range.Font.Color = 255
}

# --- 2. Relocate the _GoBack bookmark -------------------------------------
$target = "Population density (people per sq. km of land area)"
$rng2 = $d.Range(0, $d.Content.End)
$found2 = $rng2.Find.Execute($target)
if ($found2) {
    $splitPoint = $rng2.Start + "Population density".Length

    # Remove the old (hidden) _GoBack bookmark, if present.
    try {
        $existing = $d.Bookmarks("_GoBack")
        $existing.Delete()
    } catch {
    }

    # Re-create it right after "Population density".
    $bmRange = $d.Range($splitPoint, $splitPoint)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
